# Update extraction notebook and refreshed all_schools.xlsx
# - rename column M1 header (extracted variant)
# - add new column AK with the "calculated" variant of the expendable
#   net-assets-with-donor-restrictions figure

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing header in M1
$ws.Range("M1").Value = "expendable_net_assets_with_donor_restrictions_extracted"

# New header for the calculated column
$ws.Range("AK1").Value = "expendable_net_assets_with_donor_restrictions_calculated"

# Per-row calculated values (school row -> value). Rows without a
# calculated figure stay blank but still get a cell placeholder so the
# sheet's used range extends through column AK.
$akValues = @{
    2  = 777373
    3  = 43582
    4  = 146563
    5  = 5319023
    6  = 76234991
    7  = 24560
    8  = $null
    9  = $null
    10 = 4937
    11 = 64759989
    12 = $null
    13 = 1798160
    14 = 15853350
    15 = 25244510
    16 = 173476
    17 = 297801
    18 = 465238532
    19 = $null
    20 = 3771398
}

foreach ($row in $akValues.Keys) {
    $cell = $ws.Cells.Item($row, 37)
    $value = $akValues[$row]
    if ($null -eq $value) {
        # No calculated figure for this school/year: still materialize the
        # cell (so the column's used range covers it) without giving it a
        # value.
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
